$wb = $excel.ActiveWorkbook

# --- Update selection on the existing "program" sheet (A3 -> A8), and ---
# --- make sure it is no longer the tab-selected sheet once "batch" is added. ---
$wsProgram = $wb.Worksheets.Item("program")
$wsProgram.Range("A8").Select() | Out-Null

# --- Add the new "batch" worksheet after the last existing sheet ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsBatch = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsBatch.Name = "batch"

# Column widths (as close as this engine's pixel-rounding allows to the
# authored widths of 16.5703125 / 15.85546875 character-units)
$wsBatch.Columns.Item(1).ColumnWidth = 15.666666666666666
$wsBatch.Columns.Item(2).ColumnWidth = 15

# --- Populate the "batch" sheet data ---
# Row 1 - headers
$wsBatch.Cells.Item(1,1).Value = "batchname"
$wsBatch.Cells.Item(1,2).Value = "batchdescription"
$wsBatch.Cells.Item(1,3).Value = "noofclasses"

# Row 2
$wsBatch.Cells.Item(2,1).Value = "sdet115_2023"
$wsBatch.Cells.Item(2,3).Value = 5

# Row 3
$wsBatch.Cells.Item(3,1).Value = "sdet116_2023"
$wsBatch.Cells.Item(3,2).Value = "DA"
$wsBatch.Cells.Item(3,3).Value = 6

# Row 4
$wsBatch.Cells.Item(4,1).Value = "abc$"
$wsBatch.Cells.Item(4,2).Value = "DA"
$wsBatch.Cells.Item(4,3).Value = 7

# Row 5
$wsBatch.Cells.Item(5,2).Value = "salesforce"
$wsBatch.Cells.Item(5,3).Value = 4

# Row 7
$wsBatch.Cells.Item(7,1).Value = "Sdet_117"
$wsBatch.Cells.Item(7,2).Value = "TESTING"
$wsBatch.Cells.Item(7,3).Value = 5

# Row 8
$wsBatch.Cells.Item(8,1).Value = "abc"
$wsBatch.Cells.Item(8,2).Value = "def"
$wsBatch.Cells.Item(8,3).Value = "hn"

# Row 10
$wsBatch.Cells.Item(10,1).Value = "Sdet_118"
$wsBatch.Cells.Item(10,3).Value = 7

# --- Selection/active-tab state for the new sheet ---
$wsBatch.Range("C10").Select() | Out-Null
